# "sixth commit selected travelers"
#
# - ReturnDate (F2) is updated to a new date (12 October 2025 / serial 45942)
#   and picks up the "dd mmmm yyyy" custom date format (previously shared
#   with DepartureDate).
# - DepartureDate (E2) keeps its value but switches to a very similar custom
#   date format that drops the leading zero on the day ("d mmmm yyyy").
# - BoardingPlace / LandingPlace / Adults / Child (C2, D2, G2, H2) lose the
#   explicit cell style they had and fall back to the sheet's default style.
# - The saved cursor/selection moves to H8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the explicit style on the "selected travelers" cells (BoardingPlace,
# LandingPlace, Adults, Child) so they go back to the default/Normal style.
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Style = "Normal"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Style = "Normal"

# Update the return date value.
$ws.Range("F2").Value = 45942

# Re-apply the "dd mmmm yyyy" custom date format to ReturnDate (F2) first,
# then the "d mmmm yyyy" (no leading zero) custom date format to
# DepartureDate (E2) - the order matches how the two custom number formats
# end up registered in the workbook.
$ws.Range("F2").NumberFormat = "[$-14009]dd mmmm yyyy;@"
$ws.Range("E2").NumberFormat = "[$-14009]d mmmm yyyy;@"

# Match the saved selection/active cell in the worksheet view.
$ws.Range("H8").Select()
